$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '95.156.47'
$ws.Range('E2').Value = '  -2.07%  '
$ws.Range('D3').Value = '3.605.67'
$ws.Range('E3').Value = '  -2.86%  '
$ws.Range('B4').Value = 'TetherUSD'
$ws.Range('C4').Value = 'https://coinranking.com/coin/HIVsRcGKkPFtW+tetherusd-usdt'
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.998'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  -0.20%  '
$ws.Range('B5').Value = 'XRP'
$ws.Range('C5').Value = 'https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp'
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '2.27'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +19.38%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '225.33'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -5.53%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '635.94'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  -3.09%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.410'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -4.78%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '1.09'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +1.88%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.999'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +0.00%  '
$ws.Range('D11').Value = '3.595.98'
$ws.Range('E11').Value = '  -3.07%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '46.88'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +5.44%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.205'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -1.15%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.0000288'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -5.33%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '6.42'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -5.21%  '
$ws.Range('D16').Value = '4.278.72'
$ws.Range('E16').Value = '  -2.81%  '
$ws.Range('D17').Value = '94.940.05'
$ws.Range('E17').Value = '  -1.96%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '8.76'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -2.39%  '
$ws.Range('D19').Value = '3.602.39'
$ws.Range('E19').Value = '  -3.17%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '19.16'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +2.00%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '12.58'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -5.17%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.508'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -0.19%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '508.12'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -3.18%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '3.22'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -5.91%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '0.237'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +20.48%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '118.08'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +16.05%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '0.0000200'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -5.94%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '6.70'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -3.95%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '12.53'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -6.41%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '12.57'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +2.57%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '2.89'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -4.35%  '
$ws.Range('E32').Value = '  +0.08%  '
$ws.Range('B33').Value = 'Binance-PegBSC-USD'
$ws.Range('C33').Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '1.00'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -0.03%  '
$ws.Range('B34').Value = 'Cronos'
$ws.Range('C34').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.177'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -6.64%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.75'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -6.50%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '31.59'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -2.44%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.580'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -3.53%  '
$ws.Range('B38').Value = 'USDe'
$ws.Range('C38').Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '1.00'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -0.02%  '
$ws.Range('B39').Value = 'Bittensor'
$ws.Range('C39').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '591.63'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -8.90%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '8.24'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -7.13%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '6.78'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -1.08%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.484'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +8.51%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.157'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -2.55%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '38.88'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -4.70%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.0478'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +4.69%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '1.91'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -6.71%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.909'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -5.71%  '
$ws.Range('B48').Value = 'WhiteBITCoin'
$ws.Range('C48').Value = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '23.44'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -0.81%  '
$ws.Range('B49').Value = 'Cosmos'
$ws.Range('C49').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '8.50'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -0.88%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '2.18'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -4.91%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '53.46'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -0.76%  '
